$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D3").Value = 10.54
$ws.Range("E3").Value = 10.27

$ws.Range("C4").Value = 9.460000000000001
$ws.Range("E4").Value = 10.01

$ws.Range("C5").Value = 9.73
$ws.Range("D5").Value = 9.99
$ws.Range("F5").Value = 9.970000000000001

$ws.Range("E6").Value = 10.03
$ws.Range("G6").Value = 10.09
$ws.Range("H6").Value = 11.77

$ws.Range("F7").Value = 9.91
$ws.Range("H7").Value = 9.82
$ws.Range("J7").Value = 8.5

$ws.Range("F8").Value = 8.23
$ws.Range("G8").Value = 10.18

$ws.Range("G10").Value = 11.5
